# Apply cryptos list update (price/volume refresh) per commit:
# "Updated cryptos list on Tue Feb 27 18:47:48 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.045.30"
$ws.Range("E2").Value = "  +6.72%  "
$ws.Range("D3").Value = "3.239.15"
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "394.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("E7").Value = "  +5.43%  "
$ws.Range("D8").Value = "3.235.80"
$ws.Range("E8").Value = "  +3.22%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.619"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "39.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("E12").Value = "  +10.12%  "
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").Value = "3.750.97"
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "3.254.15"
$ws.Range("E17").Value = "  +3.57%  "
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").Value = "56.892.22"
$ws.Range("E20").Value = "  +6.61%  "
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  +9.02%  "
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "296.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.42%  "
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("E27").Value = "  +1.51%  "
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.27%  "
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("E34").Value = "  +2.61%  "
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.78%  "
$ws.Range("E38").Value = "  -3.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "134.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.60%  "
$ws.Range("E42").Value = "  +3.02%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "2.158.83"
$ws.Range("E48").Value = "  +4.11%  "
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +21.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.74%  "
